$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45953
$ws.Range("B2").Value = 7.22
$ws.Range("C2").Value = 1.46
$ws.Range("D2").Value = 0.22
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 2.62
$ws.Range("I2").Value = 16.49
$ws.Range("J2").Value = 23.49
$ws.Range("K2").Value = 13.83
$ws.Range("L2").Value = 10.01
$ws.Range("M2").Value = 4.67
$ws.Range("N2").Value = 4.31
$ws.Range("O2").Value = 2.03
$ws.Range("P2").Value = 0.01
$ws.Range("Q2").Value = 4.18
$ws.Range("R2").Value = 5.28
$ws.Range("S2").Value = 7.52
$ws.Range("T2").Value = 19.45
$ws.Range("U2").Value = 69.26000000000001
$ws.Range("V2").Value = 112.97
$ws.Range("W2").Value = 116.54
$ws.Range("X2").Value = 102.81
$ws.Range("Y2").Value = 76.98999999999999
$ws.Range("Z2").Value = 25.06
$ws.Range("AB2").Value = 102.33
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 114.76
$ws.Range("AE2").Value = "22h-24h"
$ws.Range("AF2").Value = 89.90000000000001
$ws.Range("AG2").Value = "0h-18h"
